# Session 2.pptx - final adjustments
$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 16: "Notebooks/" + "hyper_parameter_selection.ipynb" runs
# collapse into a single run "scripts/hyper_parameter_selection.py"
# -----------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(4)
$tr16 = $sh16.TextFrame.TextRange
$para16_9 = $tr16.Paragraphs(9, 1)
$para16_9.Text = "scripts/hyper_parameter_selection.py"

# -----------------------------------------------------------------
# Slide 20: homework slide rewritten
# -----------------------------------------------------------------
$s20 = $p.Slides.Item(20)
$sh20 = $s20.Shapes.Item(4)
$tr20 = $sh20.TextFrame.TextRange

# Paragraph 1: "Please answer the following questions:" ->
#              "Please try to find good hyper parameters"
$para1 = $tr20.Paragraphs(1, 1)
$para1.Text = "Please try to find good hyper parameters"

# Paragraph 2: "How do you evaluate trained models?" ->
#              "Please use notebook manual_model_selection.ipynb "
#              (split into 3 runs to mirror the authored formatting)
$para2 = $tr20.Paragraphs(2, 1)
$para2.Text = "Please use notebook manual_model_selection.ipynb "
$run2Full = $para2.Text
$nbName = "manual_model_selection.ipynb"
$nbStart = $para2.Start + $run2Full.IndexOf($nbName)
$nbRun = $tr20.Characters($nbStart, $nbName.Length)
$nbRun.Font.Size = 20
$trailStart = $nbStart + $nbName.Length
$trailRun = $tr20.Characters($trailStart, 1)
$trailRun.Font.Size = 20

# Paragraph 3 (old): "How do you store trained models?" is removed entirely
$para3old = $tr20.Paragraphs(3, 1)
$para3old.Delete()

# Paragraph 3 (now): "How do you use them as part of a software system?" ->
#              Change the "Model Creation" cell to find a good model
$q1 = [char]0x201C
$q2 = [char]0x201D
$para3 = $tr20.Paragraphs(3, 1)
$para3.Text = "Change the " + $q1 + "Model Creation" + $q2 + " cell to find a good model"

# New paragraph 4: "Post the cell content into this Google Form (link)"
$para3.InsertAfter([char]13 + "x")
$para4 = $tr20.Paragraphs(4, 1)
$para4.Text = "Post the cell content into this Google Form (link)"
$full4 = $para4.Text
$linkStart = $para4.Start + $full4.IndexOf("link")
$linkRun = $tr20.Characters($linkStart, 4)
$linkRun.Font.Size = 20
$actionSetting = $linkRun.ActionSettings(1)
$actionSetting.Hyperlink.Address = "https://forms.gle/"

# Paragraphs 5 and 6 (both empty) are unchanged.

# Paragraph 7: "	  Please research your answers " + (wingdings glyph) ->
#              "	  There will be a trophy for the best solution next week " + (wingdings space)
$para7 = $tr20.Paragraphs(7, 1)
$run1Len = 32
$run1 = $tr20.Characters($para7.Start, $run1Len)
$newRun1Text = [char]9 + "  There will be a trophy for the best solution next week "
$run1.Text = $newRun1Text
$wdPos = $para7.Start + $newRun1Text.Length
$wdRun = $tr20.Characters($wdPos, 1)
$wdRun.Text = " "

# Paragraph 8 (old): "             The Answers are for you. No submission needed."
# is removed; its trailing endParaRPr becomes paragraph 7's endParaRPr.
$para8old = $tr20.Paragraphs(8, 1)
$para8old.Delete()

# Remaining paragraphs (empty) are unchanged.

Write-Host "done"
